$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'double[,]' 2,100

$arr[0,0] = -7.7577058873630635
$arr[0,1] = -10.063645091669372
$arr[0,2] = -9.5360231168095364
$arr[0,3] = -8.5528597275963847
$arr[0,4] = -9.9275876512341483
$arr[0,5] = -9.2038876875467821
$arr[0,6] = -9.3210785411698112
$arr[0,7] = -7.5177794345755515
$arr[0,8] = -9.0543076084057041
$arr[0,9] = -9.2821113881334547
$arr[0,10] = -10.982219336010827
$arr[0,11] = -9.0234793453164031
$arr[0,12] = -9.7638127411866353
$arr[0,13] = -8.7580586420676525
$arr[0,14] = -9.4069421680163359
$arr[0,15] = -9.7352658048198109
$arr[0,16] = -10.845985356660158
$arr[0,17] = -10.159029609734269
$arr[0,18] = -8.4096349110532369
$arr[0,19] = -9.3012544419004684
$arr[0,20] = -9.2309956639944417
$arr[0,21] = -8.670555685738254
$arr[0,22] = -10.446912571023258
$arr[0,23] = -7.7404065594483082
$arr[0,24] = -9.845565159165103
$arr[0,25] = -9.3811575813417214
$arr[0,26] = -10.048452352767622
$arr[0,27] = -8.122167965734489
$arr[0,28] = -9.5875294108274947
$arr[0,29] = -8.8939327192725219
$arr[0,30] = -10.0076705715695
$arr[0,31] = -10.087225907005116
$arr[0,32] = -8.2783959173918227
$arr[0,33] = -9.194602125153871
$arr[0,34] = -9.117548347384913
$arr[0,35] = -8.2771072815844065
$arr[0,36] = -9.2030512227602923
$arr[0,37] = -7.3304412708172366
$arr[0,38] = -9.8794403874385051
$arr[0,39] = -10.74488654071785
$arr[0,40] = -9.0689747220739427
$arr[0,41] = -10.639002371654499
$arr[0,42] = -9.2218545872003652
$arr[0,43] = -10.344633957036761
$arr[0,44] = -8.9679344000526147
$arr[0,45] = -8.744542325049812
$arr[0,46] = -9.8515861314410653
$arr[0,47] = -10.041597223684537
$arr[0,48] = -8.9185390807639973
$arr[0,49] = -10.077192295826007
$arr[0,50] = -9.265309788365057
$arr[0,51] = -9.6102702628223629
$arr[0,52] = -8.6607877504272661
$arr[0,53] = -8.7871004663064571
$arr[0,54] = -8.8563569796938157
$arr[0,55] = -9.7923415558499229
$arr[0,56] = -9.8904699745103528
$arr[0,57] = -10.071424270590695
$arr[0,58] = -8.787066745394867
$arr[0,59] = -9.1123623215975087
$arr[0,60] = -10.416023858815752
$arr[0,61] = -8.9698379033779698
$arr[0,62] = -9.3387482881028472
$arr[0,63] = -8.154686203743271
$arr[0,64] = -10.792149665285706
$arr[0,65] = -9.8937431034280827
$arr[0,66] = -10.251628037785713
$arr[0,67] = -9.2543131727560866
$arr[0,68] = -8.5233738264615173
$arr[0,69] = -10.033454833086163
$arr[0,70] = -9.5885608046045405
$arr[0,71] = -7.8546277225847962
$arr[0,72] = -9.2030893449881823
$arr[0,73] = -8.3305247202452914
$arr[0,74] = -9.7568633313884607
$arr[0,75] = -9.7054546533571155
$arr[0,76] = -9.9992511954555408
$arr[0,77] = -9.4917639642644964
$arr[0,78] = -9.0852121284833682
$arr[0,79] = -9.2677563755747876
$arr[0,80] = -8.9039410995771533
$arr[0,81] = -10.325424925707107
$arr[0,82] = -9.329935613316108
$arr[0,83] = -9.9028142566796316
$arr[0,84] = -9.2817487398878455
$arr[0,85] = -10.180371717108072
$arr[0,86] = -9.6472737289757351
$arr[0,87] = -8.2287291907330928
$arr[0,88] = -8.675491296483349
$arr[0,89] = -11.060760012399582
$arr[0,90] = -8.9729218518285538
$arr[0,91] = -10.102727390063697
$arr[0,92] = -9.8896466434463495
$arr[0,93] = -9.4795582135907459
$arr[0,94] = -11.408333588700986
$arr[0,95] = -9.8966048705149525
$arr[0,96] = -10.470714220363611
$arr[0,97] = -9.142434265135325
$arr[0,98] = -10.001625201840547
$arr[0,99] = -9.7468989961073174

$arr[1,0] = -8.7649080954628928
$arr[1,1] = -10.101127554683577
$arr[1,2] = -9.6065074809531499
$arr[1,3] = -8.6555863365979651
$arr[1,4] = -10.055286620777581
$arr[1,5] = -9.2483885305716775
$arr[1,6] = -9.3368129449985133
$arr[1,7] = -8.6162495977730202
$arr[1,8] = -9.1188442344563647
$arr[1,9] = -9.2505989594075544
$arr[1,10] = -9.9619809098027581
$arr[1,11] = -9.0724554397263955
$arr[1,12] = -9.8592773116418151
$arr[1,13] = -8.7774938247653278
$arr[1,14] = -9.4641319040607179
$arr[1,15] = -8.8192165370324105
$arr[1,16] = -10.849729353789749
$arr[1,17] = -10.150341818182724
$arr[1,18] = -9.3608273763299916
$arr[1,19] = -9.2449310779528258
$arr[1,20] = -9.140754591173522
$arr[1,21] = -9.790094504143811
$arr[1,22] = -10.447930619701479
$arr[1,23] = -8.7643421522458969
$arr[1,24] = -9.7663691502483321
$arr[1,25] = -9.4457307214691255
$arr[1,26] = -10.195264248359486
$arr[1,27] = -8.1274899409580357
$arr[1,28] = -9.6554594800256783
$arr[1,29] = -8.8702148257972624
$arr[1,30] = -10.008066462338666
$arr[1,31] = -10.103315477101784
$arr[1,32] = -9.2481974711406529
$arr[1,33] = -8.2588458346059266
$arr[1,34] = -9.0391628385990508
$arr[1,35] = -9.2074192600240323
$arr[1,36] = -9.2858402418980148
$arr[1,37] = -8.330784959844129
$arr[1,38] = -9.8705843388861734
$arr[1,39] = -9.66541414188411
$arr[1,40] = -8.0380909506651363
$arr[1,41] = -10.615247025271023
$arr[1,42] = -8.1763017835157044
$arr[1,43] = -10.363760388215322
$arr[1,44] = -8.1101215841122443
$arr[1,45] = -9.9122734216838548
$arr[1,46] = -9.8104028211489638
$arr[1,47] = -10.048589785497807
$arr[1,48] = -8.8625231153301804
$arr[1,49] = -9.1767943957817995
$arr[1,50] = -9.2862194982771946
$arr[1,51] = -8.5824881203179899
$arr[1,52] = -8.7302189585928627
$arr[1,53] = -8.8061109782733951
$arr[1,54] = -8.8694109391593106
$arr[1,55] = -10.8302375228398
$arr[1,56] = -9.8783101434024161
$arr[1,57] = -9.0695066164254428
$arr[1,58] = -9.7868542604221993
$arr[1,59] = -9.2305606771165394
$arr[1,60] = -10.458590686242013
$arr[1,61] = -9.9862820775240024
$arr[1,62] = -9.3530588251116118
$arr[1,63] = -8.1738082418064444
$arr[1,64] = -10.794887920647152
$arr[1,65] = -8.9106933934330073
$arr[1,66] = -10.348533357540459
$arr[1,67] = -9.2725593022190491
$arr[1,68] = -9.5583292284498409
$arr[1,69] = -10.106664530036406
$arr[1,70] = -9.7095014188686495
$arr[1,71] = -8.8400151414748276
$arr[1,72] = -9.2322994290598714
$arr[1,73] = -9.4400633470640347
$arr[1,74] = -9.7856629491406668
$arr[1,75] = -9.7309457232618026
$arr[1,76] = -10.14668177060787
$arr[1,77] = -9.6216020164790503
$arr[1,78] = -9.0934095696409649
$arr[1,79] = -9.284069706843713
$arr[1,80] = -10.045404078307341
$arr[1,81] = -10.294311323493384
$arr[1,82] = -9.359304550371581
$arr[1,83] = -9.8298226895599505
$arr[1,84] = -8.2938190540745289
$arr[1,85] = -10.153429898227012
$arr[1,86] = -9.6498275692371909
$arr[1,87] = -9.2520228546828349
$arr[1,88] = -8.7481668305747675
$arr[1,89] = -11.04730585335742
$arr[1,90] = -10.005167504734979
$arr[1,91] = -9.1954893511310214
$arr[1,92] = -9.9825055368004598
$arr[1,93] = -9.5724941423796341
$arr[1,94] = -11.304821413151272
$arr[1,95] = -9.8958977312182306
$arr[1,96] = -10.441499651450755
$arr[1,97] = -8.117141530240179
$arr[1,98] = -9.9498877766991356
$arr[1,99] = -9.8894092350000111

$ws.Range("A1:CV2").Value = $arr

Write-Output "Done"